$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: the login username field becomes the twitter-id field, and the
# Pass/Fail column (G) gets a "Pass" result.
$ws.Range("B2").Value = "loginData.twitterID"

# Give the new Pass/Fail cells the same look (style) as the rest of column
# G / F by copying the format from an already-styled neighbour, then fill
# in the values.
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "Pass"

$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "Pass"

$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Range("F6").Copy()
$ws.Range("G6").PasteSpecial(-4122)

$ws.Range("F8").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("F10").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("G10").Value = "Pass"

$excel.CutCopyMode = 0
